$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18 and 19 swap coin identity (Dai <-> Uniswap) in addition to the
# price/volume refresh applied to every data row below.
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"

# Price (D) and Volume(1h) (E) refresh for every data row (2-51).
# D values are written with a leading apostrophe so Excel keeps them as
# literal text (several look numeric, e.g. "0.9991", "1.000") and the
# style is immediately reset to Normal so no stray text-format style
# sticks to the cell.
$ws.Range("D2").Value = "'30.325.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "'1.870.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("D4").Value = "'0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'235.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.49%  "
$ws.Range("D6").Value = "'0.9987"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.4676"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("D8").Value = "'0.2846"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.06566"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'20.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.58%  "
$ws.Range("D11").Value = "'0.07888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").Value = "'97.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "'1.866.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "'5.161"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "'0.6772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'283.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'30.315.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "'5.534"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.58%  "
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'12.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("D21").Value = "'2.109.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'0.000007284"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").Value = "'0.9991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("D24").Value = "'6.195"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "'9.328"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'165.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "'19.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'1.916"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.80%  "
$ws.Range("D29").Value = "'1.355"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "'0.09688"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("D31").Value = "'4.441"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "'1.472"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").Value = "'4.113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.97%  "
$ws.Range("D34").Value = "'0.04720"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "'1.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").Value = "'0.7055"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").Value = "'2.715"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").Value = "'0.01863"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").Value = "'6.386"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.48%  "
$ws.Range("D40").Value = "'2.534"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("D41").Value = "'73.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "'1.946"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'0.8498"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "'0.4194"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "'104.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.55%  "
$ws.Range("D46").Value = "'0.9990"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'7.225"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.96%  "
$ws.Range("D48").Value = "'9.317"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'946.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("D50").Value = "'34.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "'0.1136"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.72%  "
